# Validate Supervisor Ticket List Page
#
# 1) Update the LoginCredentials sheet's sample row (new login id / password /
#    user type).
# 2) Add a new "TicketId" worksheet (after LoginCredentials) holding a header
#    + a ticket id that must keep its leading zero (quote-prefixed text).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- LoginCredentials updates ---------------------------------------------
$ws1.Range("A2").Value = 2390495
$ws1.Range("B2").Value = "June@123$"
$ws1.Range("C2").Value = "BS"

$ws1.Range("C6").Select() | Out-Null

# --- New TicketId sheet -----------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TicketId"

$ws2.Range("A1").Value = "ticketid"
# Leading apostrophe forces text storage (quote-prefixed) so the leading
# zero in the ticket id is preserved instead of being parsed as a number.
$ws2.Range("A2").Value = "'080720000457"

$ws2.Columns.Item(1).ColumnWidth = 12.43

$ws2.Range("E4").Select() | Out-Null
